# Regenerate save_data: recalc the "K" column (G) values with updated figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 4
    4  = 5
    5  = 6
    6  = 6
    7  = 5
    8  = 6
    9  = 5
    10 = 6
    12 = 2
    13 = 1
    14 = 0
    15 = 2
    16 = 6
    17 = 4
    19 = 1
    20 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
